# Daily attendance processing - swap the order of "System" and the email
# address in the "Recorded By" column (G) wherever both are present,
# combined as "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
